$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" sheet: insert a new row 2 for the "2022-Q4" summary record,
#    pushing the existing quarters down by one row and renumbering the
#    running index in column A.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Range("B2:D2").ClearFormats()
$summary.Range("A3").Copy($summary.Range("A2"))

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.18

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

# ---------------------------------------------------------------------
# 2. Add the new "2022-Q4" sheet right after "总计", using the existing
#    "2021-Q2" sheet as a formatting template (same header layout,
#    column widths, page margins, etc.) then overwrite its data.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q2")
$template.Copy($null, $summary)

$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q4"

# Only one data row is needed for this quarter - drop the extra template rows.
$newSheet.Range("A3:A5").EntireRow.Delete()

# "基金规模" replaces "基金金额" as the column D header for this quarter.
$newSheet.Range("D1").Value = "基金规模"

$newSheet.Range("A2").Value = 0

# Force these as text (not numbers) so values like "003956" / "85.80"
# keep their leading/trailing zeros, then drop the resulting "@" style
# so the cells stay unformatted like the rest of the sheet.
$newSheet.Range("B2:G2").NumberFormat = "@"
$newSheet.Range("B2").Value = "003956"
$newSheet.Range("C2").Value = "南方产业智选股票"
$newSheet.Range("D2").Value = "3.60"
$newSheet.Range("E2").Value = "85.80"
$newSheet.Range("F2").Value = "4.90"
$newSheet.Range("G2").Value = "0.1764"
$newSheet.Range("B2:G2").ClearFormats()

$newSheet.Range("H2").Value = 3

# ---------------------------------------------------------------------
# 3. Restore the original tab selection (the last sheet, "2020-Q4",
#    was the selected tab before this edit).
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
